$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.761.16"
$ws.Range("D2").NumberFormat = "General"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -3.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.799.31"
$ws.Range("D3").NumberFormat = "General"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +1.51%  "

$ws.Range("E4").Value = "  +0.18%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.99"
$ws.Range("D5").NumberFormat = "General"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -3.25%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "173.45"
$ws.Range("D6").NumberFormat = "General"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -3.77%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.798.76"
$ws.Range("D7").NumberFormat = "General"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  +1.47%  "

$ws.Range("E8").Value = "  +0.10%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.526"
$ws.Range("D9").NumberFormat = "General"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.98%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.159"
$ws.Range("D10").NumberFormat = "General"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -4.03%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.26"
$ws.Range("D11").NumberFormat = "General"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -5.04%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.462"
$ws.Range("D12").NumberFormat = "General"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -4.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "37.94"
$ws.Range("D13").NumberFormat = "General"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -5.09%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.0000244"
$ws.Range("D14").NumberFormat = "General"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -4.31%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.446.21"
$ws.Range("D15").NumberFormat = "General"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +1.84%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.807.09"
$ws.Range("D16").NumberFormat = "General"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +1.61%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "67.939.52"
$ws.Range("D17").NumberFormat = "General"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  -2.76%  "

$ws.Range("E18").Value = "  -4.53%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.13"
$ws.Range("D19").NumberFormat = "General"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -5.56%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "16.23"
$ws.Range("D20").NumberFormat = "General"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  -1.83%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "488.87"
$ws.Range("D21").NumberFormat = "General"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -2.96%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "9.18"
$ws.Range("D22").NumberFormat = "General"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.14%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.728"
$ws.Range("D23").NumberFormat = "General"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  +0.59%  "

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.47"
$ws.Range("D24").NumberFormat = "General"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -2.36%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.37"
$ws.Range("D25").NumberFormat = "General"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -9.47%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0000137"
$ws.Range("D26").NumberFormat = "General"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +1.78%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.29"
$ws.Range("D27").NumberFormat = "General"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -5.65%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.17"
$ws.Range("D28").NumberFormat = "General"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -10.30%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.00"
$ws.Range("D29").NumberFormat = "General"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("E30").Value = "  -0.18%  "

$ws.Range("E31").Value = "  -2.57%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "32.67"
$ws.Range("D32").NumberFormat = "General"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +6.54%  "

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "7.69"
$ws.Range("D33").NumberFormat = "General"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -3.76%  "

$ws.Range("E34").Value = "  -4.59%  "

$ws.Range("E35").Value = "  +0.03%  "

$ws.Range("E36").Value = "  -4.71%  "

$ws.Range("E37").Value = "  -1.96%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "5.76"
$ws.Range("D38").NumberFormat = "General"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -6.31%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.324"
$ws.Range("D39").NumberFormat = "General"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -7.34%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "445.63"
$ws.Range("D40").NumberFormat = "General"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +1.50%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "48.87"
$ws.Range("D41").NumberFormat = "General"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -1.61%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "1.98"
$ws.Range("D42").NumberFormat = "General"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -4.23%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.88"
$ws.Range("D43").NumberFormat = "General"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -6.40%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "8.24"
$ws.Range("D44").NumberFormat = "General"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -4.38%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "41.27"
$ws.Range("D45").NumberFormat = "General"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -9.97%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "2.829.38"
$ws.Range("D46").NumberFormat = "General"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -4.85%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "138.55"
$ws.Range("D48").NumberFormat = "General"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -0.50%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0350"
$ws.Range("D49").NumberFormat = "General"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -3.96%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "26.03"
$ws.Range("D50").NumberFormat = "General"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -4.75%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.30"
$ws.Range("D51").NumberFormat = "General"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -6.96%  "

